# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# --- 1. Update the "总计" sheet: insert a new top data row for 2022-Q1 ---
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Whole-row insert copies the header's bold/centered look onto the new row;
# strip that back off and re-apply the plain row-index style (from A3, which
# still carries the original formatting) to the new A2 cell only.
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.12

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

# --- 2. Insert the new "2022-Q1" sheet right before "总计" ---------------
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Match the look of the other quarter sheets: bordered/bold header style and
# the same row-index style used for column A.
$wb.Worksheets.Item("2021-Q4").Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$wb.Worksheets.Item("2021-Q4").Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'002236"
$newSheet.Range("C2").Value = "大成中证360互联网+大数据100指数A"
$newSheet.Range("D2").Value = "'5.67"
$newSheet.Range("E2").Value = "'93.32"
$newSheet.Range("F2").Value = "'1.21"
$newSheet.Range("G2").Value = "'0.0686"
$newSheet.Range("H2").Value = 1

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'003359"
$newSheet.Range("C3").Value = "大成中证360互联网+大数据100指数C"
$newSheet.Range("D3").Value = "'4.08"
$newSheet.Range("E3").Value = "'93.32"
$newSheet.Range("F3").Value = "'1.21"
$newSheet.Range("G3").Value = "'0.0494"
$newSheet.Range("H3").Value = 1

# Restore the original active sheet (2021-Q1) so the workbook-level view
# state is left untouched by this edit.
$wb.Worksheets.Item("2021-Q1").Activate()
